$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 375, pushing existing rows 375-401 down to 376-402.
$ws.Rows.Item(375).Insert()

# Populate the newly inserted row 375 with the new weekly record
# (same Mercado/Región/Categoría context as the surrounding rows).
$ws.Cells.Item(375, 1).Value = 5
$ws.Cells.Item(375, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(375, 3).Value = "Maule"
$ws.Cells.Item(375, 4).Value = 45021
$ws.Cells.Item(375, 5).Value = 7
$ws.Cells.Item(375, 6).Value = 100112008
$ws.Cells.Item(375, 7).Value = "Coliflor"
$ws.Cells.Item(375, 8).Value = "Sin especificar"
$ws.Cells.Item(375, 9).Value = "Primera"
$ws.Cells.Item(375, 10).Value = 3000
$ws.Cells.Item(375, 11).Value = 800
$ws.Cells.Item(375, 12).Value = 800
$ws.Cells.Item(375, 13).Value = 800
$ws.Cells.Item(375, 14).Value = "`$/unidad"
$ws.Cells.Item(375, 15).Value = "Región del Maule"
$ws.Cells.Item(375, 16).Value = 800
$ws.Cells.Item(375, 17).Value = 1
$ws.Cells.Item(375, 18).Value = "Hortaliza"

# Match the date style used by the rest of column D (numFmt YYYY-MM-DD HH:MM:SS).
$ws.Cells.Item(375, 4).NumberFormat = $ws.Cells.Item(376, 4).NumberFormat
